# Update LinkedIn job posting form: append a new data row (row 7) for the
# "Landing Gear System Engineer" posting, reusing the job-description text
# that's already used by rows 3, 4 and 6 so the shared-string table isn't
# duplicated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy/PasteSpecial (instead of a direct .Value assignment) reuses the
# existing shared-string entry and avoids an unwanted row-height autofit
# that a plain multi-line .Value write would trigger.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial()

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 3
